{"js": "// Remove \"\uff0c\u8c03\u8bd5\u7ec6\u8282\" from the last progress-note paragraph, turning:\n//   \"...\u6574\u4f53\u529f\u80fd\u8fd0\u884c\u4e5f\u6b63\u5e38\uff0c\u8c03\u8bd5\u7ec6\u8282\uff0c\u7136\u540e\u628a\u534f\u8bae\u6539\u4e3aFTP\u534f\u8bae\u3002\"\n// into:\n//   \"...\u6574\u4f53\u529f\u80fd\u8fd0\u884c\u4e5f\u6b63\u5e38\uff0c\u7136\u540e\u628a\u534f\u8bae\u6539\u4e3aFTP\u534f\u8bae\u3002\"\n//\n// Word re-anchors its \"_GoBack\" (last-edit) bookmark to the spot of the\n// edit, so after the deletion the bookmark sits between the two remaining\n// halves of the sentence instead of at the end of the paragraph. We\n// reproduce that by moving the bookmark along with the text edit.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// Locate the exact span that was removed. Anchoring the search on the\n// leading \"\uff0c\" (rather than the trailing one) keeps the surviving comma on\n// the correct side of the cut, matching \"...\u4e5f\u6b63\u5e38\" + \"\uff0c\u7136\u540e...\".\nconst results = body.search(\"\uff0c\u8c03\u8bd5\u7ec6\u8282\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Caret marking where the remaining text will be rejoined - this is\n  // where the \"_GoBack\" bookmark needs to live once the old one is gone.\n  const splitPoint = target.getRange(\"Start\");\n\n  // The existing \"_GoBack\" bookmark currently sits at the end of the\n  // paragraph; drop it before we touch the text.\n  doc.deleteBookmark(\"_GoBack\");\n\n  // Remove \"\uff0c\u8c03\u8bd5\u7ec6\u8282\" from the sentence.\n  target.delete();\n  await context.sync();\n\n  // Re-create \"_GoBack\" at the edit point (between the two remaining runs).\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Remove \"\uff0c\u8c03\u8bd5\u7ec6\u8282\" from the last progress-note paragraph, turning:\n#   \"...\u6574\u4f53\u529f\u80fd\u8fd0\u884c\u4e5f\u6b63\u5e38\uff0c\u8c03\u8bd5\u7ec6\u8282\uff0c\u7136\u540e\u628a\u534f\u8bae\u6539\u4e3aFTP\u534f\u8bae\u3002\"\n# into:\n#   \"...\u6574\u4f53\u529f\u80fd\u8fd0\u884c\u4e5f\u6b63\u5e38\uff0c\u7136\u540e\u628a\u534f\u8bae\u6539\u4e3aFTP\u534f\u8bae\u3002\"\n#\n# Word re-anchors its \"_GoBack\" (last-edit) bookmark to the spot of the\n# edit, so after the deletion the bookmark sits between the two remaining\n# halves of the sentence instead of at the end of the paragraph. We\n# reproduce that by moving the bookmark along with the text edit.\n\n$d = $word.ActiveDocument\n\n# The existing \"_GoBack\" bookmark currently sits at the end of the\n# paragraph; drop it before we touch the text (Word will recreate it at\n# the new edit point).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the exact span that was removed. Anchoring the search on the\n# leading \"\uff0c\" (rather than the trailing one) keeps the surviving comma on\n# the correct side of the cut, matching \"...\u4e5f\u6b63\u5e38\" + \"\uff0c\u7136\u540e...\".\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\uff0c\u8c03\u8bd5\u7ec6\u8282\")\n\nif ($found) {\n    # Caret marking where the remaining text will be rejoined - this is\n    # where the \"_GoBack\" bookmark needs to live once it's re-created.\n    $splitRange = $d.Range($rng.Start, $rng.Start)\n\n    # Remove \"\uff0c\u8c03\u8bd5\u7ec6\u8282\" from the sentence.\n    $rng.Text = \"\"\n\n    # Re-create \"_GoBack\" at the edit point (between the two remaining runs).\n    $d.Bookmarks.Add(\"_GoBack\", $splitRange)\n}\n"}
